# Generate Report for Handback
# Update the "Correspond Handoff Datetime" / "Correspond Handback DateTime"
# timestamps on the zh-cn and de-de worksheets to reflect the newly
# generated handback report timings.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E2").Value = "2016-03-21 23:04:57"
$wsZhCn.Range("H2").Value = "2016-03-21 23:05:24"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E2").Value = "2016-03-21 23:05:01"
$wsDeDe.Range("H2").Value = "2016-03-21 23:05:30"
